$d = $word.ActiveDocument

$d.Content.Find.Execute("817×9=7353", $true, $false, $false, $false, $false, $true, 1, $false, "988×5=4940", 2) | Out-Null
$d.Content.Find.Execute("383×3=1149", $true, $false, $false, $false, $false, $true, 1, $false, "824×2=1648", 2) | Out-Null
$d.Content.Find.Execute("145×7=1015", $true, $false, $false, $false, $false, $true, 1, $false, "623×9=5607", 2) | Out-Null
$d.Content.Find.Execute("563×6=3378", $true, $false, $false, $false, $false, $true, 1, $false, "567×6=3402", 2) | Out-Null
$d.Content.Find.Execute("157×8=1256", $true, $false, $false, $false, $false, $true, 1, $false, "885×8=7080", 2) | Out-Null
$d.Content.Find.Execute("482×4=1928", $true, $false, $false, $false, $false, $true, 1, $false, "839×3=2517", 2) | Out-Null
$d.Content.Find.Execute("404×5=2020", $true, $false, $false, $false, $false, $true, 1, $false, "696×4=2784", 2) | Out-Null
$d.Content.Find.Execute("767×4=3068", $true, $false, $false, $false, $false, $true, 1, $false, "366×3=1098", 2) | Out-Null
$d.Content.Find.Execute("601×3=1803", $true, $false, $false, $false, $false, $true, 1, $false, "671×8=5368", 2) | Out-Null
$d.Content.Find.Execute("561×9=5049", $true, $false, $false, $false, $false, $true, 1, $false, "766×8=6128", 2) | Out-Null
$d.Content.Find.Execute("856×3=2568", $true, $false, $false, $false, $false, $true, 1, $false, "418×2=836", 2) | Out-Null
$d.Content.Find.Execute("446×3=1338", $true, $false, $false, $false, $false, $true, 1, $false, "761×7=5327", 2) | Out-Null
$d.Content.Find.Execute("662×7=4634", $true, $false, $false, $false, $false, $true, 1, $false, "263×7=1841", 2) | Out-Null
$d.Content.Find.Execute("649×5=3245", $true, $false, $false, $false, $false, $true, 1, $false, "619×7=4333", 2) | Out-Null
$d.Content.Find.Execute("348×3=1044", $true, $false, $false, $false, $false, $true, 1, $false, "526×7=3682", 2) | Out-Null
$d.Content.Find.Execute("575×3=1725", $true, $false, $false, $false, $false, $true, 1, $false, "850×7=5950", 2) | Out-Null
$d.Content.Find.Execute("687×3=2061", $true, $false, $false, $false, $false, $true, 1, $false, "247×9=2223", 2) | Out-Null
$d.Content.Find.Execute("720×5=3600", $true, $false, $false, $false, $false, $true, 1, $false, "450×9=4050", 2) | Out-Null
$d.Content.Find.Execute("122×8=976", $true, $false, $false, $false, $false, $true, 1, $false, "869×9=7821", 2) | Out-Null
$d.Content.Find.Execute("695×3=2085", $true, $false, $false, $false, $false, $true, 1, $false, "873×4=3492", 2) | Out-Null
$d.Content.Find.Execute("990×5=4950", $true, $false, $false, $false, $false, $true, 1, $false, "908×7=6356", 2) | Out-Null
$d.Content.Find.Execute("420×5=2100", $true, $false, $false, $false, $false, $true, 1, $false, "194×7=1358", 2) | Out-Null
$d.Content.Find.Execute("775×9=6975", $true, $false, $false, $false, $false, $true, 1, $false, "194×2=388", 2) | Out-Null
$d.Content.Find.Execute("480×5=2400", $true, $false, $false, $false, $false, $true, 1, $false, "122×4=488", 2) | Out-Null
$d.Content.Find.Execute("918×2=1836", $true, $false, $false, $false, $false, $true, 1, $false, "743×2=1486", 2) | Out-Null
